$d = $word.ActiveDocument

# --- 1. Mark "Trainingspunkte müssen erklärt werden" as explained: strike it
#        through (stamps both the paragraph-mark rPr and the run rPr, same
#        as Word does when you select the line and hit Strikethrough). ---
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Trainingspunkte müssen erklärt werden") {
        $p.Range.Font.StrikeThrough = $true
        break
    }
}

# --- 2. Move the "_GoBack" bookmark from the end of "Punkte entfernen für
#        Auf- und Abwärmübungen" to the end of "Who you are Mihail". ---

# Locate the "Who you are Mihail" paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Who you are Mihail") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a throw-away sentinel character right after the run text so the
    # bookmark's insertion point is no longer the very last offset of the
    # paragraph (that boundary makes bookmark insertion land before the run
    # instead of after it). We delete the sentinel right after.
    $ip = $target.Range.Duplicate
    $ip.MoveEnd(1, -1)
    $ip.Collapse(0)
    $ip.InsertAfter([char]1)

    $p2 = $target.Range.Duplicate
    $p2.MoveEnd(1, -1)
    $bmStart = $p2.End - 1
    $bmRange = $d.Range($bmStart, $bmStart)

    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $sentinel = $d.Range($bmStart, $bmStart + 1)
    $sentinel.Delete()
}
